{"js": "// Replace each \"A\u00f7B=C, D\" answer string in the practice-sheet table with\n// its updated value, per the commit's regenerated output.\nconst replacements = [\n  [\"127\u00f77=18, 1\", \"780\u00f72=390, 0\"],\n  [\"437\u00f77=62, 3\", \"477\u00f72=238, 1\"],\n  [\"335\u00f79=37, 2\", \"415\u00f78=51, 7\"],\n  [\"808\u00f77=115, 3\", \"271\u00f72=135, 1\"],\n  [\"734\u00f76=122, 2\", \"529\u00f72=264, 1\"],\n  [\"469\u00f74=117, 1\", \"201\u00f79=22, 3\"],\n  [\"784\u00f75=156, 4\", \"310\u00f77=44, 2\"],\n  [\"208\u00f75=41, 3\", \"911\u00f78=113, 7\"],\n  [\"834\u00f77=119, 1\", \"413\u00f78=51, 5\"],\n  [\"697\u00f77=99, 4\", \"947\u00f73=315, 2\"],\n  [\"897\u00f74=224, 1\", \"848\u00f79=94, 2\"],\n  [\"767\u00f74=191, 3\", \"682\u00f74=170, 2\"],\n  [\"113\u00f79=12, 5\", \"863\u00f77=123, 2\"],\n  [\"380\u00f78=47, 4\", \"970\u00f75=194, 0\"],\n  [\"498\u00f77=71, 1\", \"668\u00f76=111, 2\"],\n  [\"888\u00f76=148, 0\", \"168\u00f73=56, 0\"],\n  [\"489\u00f72=244, 1\", \"832\u00f76=138, 4\"],\n  [\"302\u00f74=75, 2\", \"727\u00f76=121, 1\"],\n  [\"202\u00f74=50, 2\", \"149\u00f76=24, 5\"],\n  [\"924\u00f75=184, 4\", \"337\u00f74=84, 1\"],\n  [\"120\u00f79=13, 3\", \"384\u00f75=76, 4\"],\n  [\"881\u00f76=146, 5\", \"912\u00f72=456, 0\"],\n  [\"893\u00f72=446, 1\", \"530\u00f77=75, 5\"],\n  [\"677\u00f78=84, 5\", \"907\u00f73=302, 1\"],\n  [\"498\u00f76=83, 0\", \"855\u00f77=122, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"A\u00f7B=C, D\" answer string in the practice-sheet table with\n# its updated value, per the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"127\u00f77=18, 1\", \"780\u00f72=390, 0\"),\n    @(\"437\u00f77=62, 3\", \"477\u00f72=238, 1\"),\n    @(\"335\u00f79=37, 2\", \"415\u00f78=51, 7\"),\n    @(\"808\u00f77=115, 3\", \"271\u00f72=135, 1\"),\n    @(\"734\u00f76=122, 2\", \"529\u00f72=264, 1\"),\n    @(\"469\u00f74=117, 1\", \"201\u00f79=22, 3\"),\n    @(\"784\u00f75=156, 4\", \"310\u00f77=44, 2\"),\n    @(\"208\u00f75=41, 3\", \"911\u00f78=113, 7\"),\n    @(\"834\u00f77=119, 1\", \"413\u00f78=51, 5\"),\n    @(\"697\u00f77=99, 4\", \"947\u00f73=315, 2\"),\n    @(\"897\u00f74=224, 1\", \"848\u00f79=94, 2\"),\n    @(\"767\u00f74=191, 3\", \"682\u00f74=170, 2\"),\n    @(\"113\u00f79=12, 5\", \"863\u00f77=123, 2\"),\n    @(\"380\u00f78=47, 4\", \"970\u00f75=194, 0\"),\n    @(\"498\u00f77=71, 1\", \"668\u00f76=111, 2\"),\n    @(\"888\u00f76=148, 0\", \"168\u00f73=56, 0\"),\n    @(\"489\u00f72=244, 1\", \"832\u00f76=138, 4\"),\n    @(\"302\u00f74=75, 2\", \"727\u00f76=121, 1\"),\n    @(\"202\u00f74=50, 2\", \"149\u00f76=24, 5\"),\n    @(\"924\u00f75=184, 4\", \"337\u00f74=84, 1\"),\n    @(\"120\u00f79=13, 3\", \"384\u00f75=76, 4\"),\n    @(\"881\u00f76=146, 5\", \"912\u00f72=456, 0\"),\n    @(\"893\u00f72=446, 1\", \"530\u00f77=75, 5\"),\n    @(\"677\u00f78=84, 5\", \"907\u00f73=302, 1\"),\n    @(\"498\u00f76=83, 0\", \"855\u00f77=122, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
